$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data: Language (A2:A23), Value (B2:B23)
$data = @(
    [PSCustomObject]@{Name="Arabic"; Value=4.620048191165291},
    [PSCustomObject]@{Name="Bengali"; Value=0.4451942296638969},
    [PSCustomObject]@{Name="Chinese"; Value=9.457340532913307},
    [PSCustomObject]@{Name="Dutch"; Value=1.609898565311415},
    [PSCustomObject]@{Name="English"; Value=26.74127979306701},
    [PSCustomObject]@{Name="French"; Value=3.451316456163823},
    [PSCustomObject]@{Name="German"; Value=5.620818902120796},
    [PSCustomObject]@{Name="Italian"; Value=3.195826498766609},
    [PSCustomObject]@{Name="Japanese"; Value=6.565396164626468},
    [PSCustomObject]@{Name="Korean"; Value=1.676149646860095},
    [PSCustomObject]@{Name="Malay-Indonesian"; Value=2.615870126511768},
    [PSCustomObject]@{Name="Persian"; Value=1.518461750752512},
    [PSCustomObject]@{Name="Polish"; Value=0.8707222916192066},
    [PSCustomObject]@{Name="Portuguese"; Value=3.660639813232977},
    [PSCustomObject]@{Name="Russian"; Value=3.303284177572094},
    [PSCustomObject]@{Name="Spanish"; Value=7.604050506195568},
    [PSCustomObject]@{Name="Swedish"; Value=0.5090549482866712},
    [PSCustomObject]@{Name="Thai"; Value=0.9482424856587278},
    [PSCustomObject]@{Name="Turkish"; Value=1.321307514033423},
    [PSCustomObject]@{Name="Urdu"; Value=0.8029500700267684},
    [PSCustomObject]@{Name="Uzbek"; Value=0.134023303141118},
    [PSCustomObject]@{Name="Vietnamese"; Value=0.4254139764520277}
)

# Sort descending by value, keep the top 20 rows (drop the two smallest)
$sorted = $data | Sort-Object -Property Value -Descending
$top20 = $sorted | Select-Object -First 20

# Write sorted data back starting at row 2
$r = 2
foreach ($item in $top20) {
    $ws.Cells.Item($r, 1).Value = $item.Name
    $ws.Cells.Item($r, 2).Value = $item.Value
    $r = $r + 1
}

# Remove the now-unused rows 22 and 23 (previously held the two lowest entries)
$ws.Range("A22:B23").Delete() | Out-Null
